# The workbook tracks weekly price observations for "Repollo" (cabbage) at
# the "Terminal Hortofrutícola Agro Chillán" market. A new weekly
# observation is inserted at row 26, which pushes every subsequent
# observation (previously rows 26-131) down by one row, ending with a new
# last row 132 (holding what used to be row 131's data). The sheet's
# dimension grows from A1:R131 to A1:R132 automatically as part of the
# insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26. This shifts rows 26:131 down to 27:132
# and carries the existing row formatting (e.g. the date style on column D)
# down with them, exactly like Excel's own "Insert Sheet Rows" command.
$ws.Rows.Item(26).Insert()

# Columns that stay the same for this new observation (market/region/
# product/variety/quality/unit/origin/classification) - copy them from the
# row directly below (row 27), which still holds the values that used to
# live in row 26 before the insert.
$ws.Range("A26").Value = $ws.Range("A27").Value2
$ws.Range("B26").Value = $ws.Range("B27").Value2
$ws.Range("C26").Value = $ws.Range("C27").Value2
$ws.Range("E26").Value = $ws.Range("E27").Value2
$ws.Range("F26").Value = $ws.Range("F27").Value2
$ws.Range("G26").Value = $ws.Range("G27").Value2
$ws.Range("H26").Value = $ws.Range("H27").Value2
$ws.Range("I26").Value = $ws.Range("I27").Value2
$ws.Range("N26").Value = $ws.Range("N27").Value2
$ws.Range("O26").Value = $ws.Range("O27").Value2
$ws.Range("Q26").Value = $ws.Range("Q27").Value2
$ws.Range("R26").Value = $ws.Range("R27").Value2

# New values reported for this observation: date, volume, min/max/avg
# price and the $/Kg column (which mirrors the average price here).
$ws.Range("D26").Value = 44487
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 600
$ws.Range("L26").Value = 700
$ws.Range("M26").Value = 650
$ws.Range("P26").Value = 650
